# TC05_CDS_phs001524_SampleType_Blood.xlsx
# "Added CDS All studies testcase"
#
# The SamplesTab query (cell B3) is rewritten to drop the two extra
# columns (Tumor / Analyte Type) that used to be pulled from the sample
# table, keeping just Sample ID / Participant ID / Study Name / Accession.
# The FilesTab query in B4 is untouched - it only moves around inside the
# shared-string table (an OOXML-save side effect), which Excel/the engine
# handles on its own once B3's text actually changes.
# The active selection also moves from C4 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSamplesTabQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001524' AND smp.sample_type = 'Blood'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesTabQuery

# Move the active cell/selection to B3 (was C4) and scroll the viewport so
# row 3 is at the top (best-effort - not every host persists scroll state).
$ws.Range("B3").Select()
try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
